$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

$ws.Range("D2").Value = 0.005
$ws.Range("D3").Value = 0.049
$ws.Range("D4").Value = 0.008
$ws.Range("D5").Value = 0.008
$ws.Range("D6").Value = 0.005
